$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header labels ---
# The fill order below reproduces the exact shared-string insertion order
# recorded in the target workbook: the ".25."/".75." sub-headers are written
# first (P..T then V..Z), and the "25 Quantile"/"75 Quantile" group headers
# (O1/U1) are written last.
$ws.Range("P1").Value = "specificity.25."
$ws.Range("Q1").Value = "sensitivity.25."
$ws.Range("R1").Value = "accuracy.25."
$ws.Range("S1").Value = "threshold.25."
$ws.Range("T1").Value = "AUC.25."
$ws.Range("V1").Value = "specificity.75."
$ws.Range("W1").Value = "sensitivity.75."
$ws.Range("X1").Value = "accuracy.75."
$ws.Range("Y1").Value = "threshold.75."
$ws.Range("Z1").Value = "AUC.75."
$ws.Range("O1").Value = "25 Quantile"
$ws.Range("U1").Value = "75 Quantile"

# Group header cells are bold, matching the existing "3 Mean" / "3 Median" headers
$ws.Range("O1").Font.Bold = $true
$ws.Range("U1").Font.Bold = $true

# --- Data rows: specificity / sensitivity / accuracy / threshold / AUC ---
# for the 25th percentile (P:T) and 75th percentile (V:Z) summaries
# Row 2
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 97
$ws.Range("S2").Value = 50
$ws.Range("T2").Value = 78
$ws.Range("V2").Value = 100
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 97
$ws.Range("Y2").Value = 50
$ws.Range("Z2").Value = 80

# Row 3
$ws.Range("P3").Value = 72
$ws.Range("Q3").Value = 67
$ws.Range("R3").Value = 72
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 78
$ws.Range("V3").Value = 78
$ws.Range("W3").Value = 74
$ws.Range("X3").Value = 78
$ws.Range("Y3").Value = 3
$ws.Range("Z3").Value = 80

# Row 4
$ws.Range("P4").Value = 71
$ws.Range("Q4").Value = 66
$ws.Range("R4").Value = 71
$ws.Range("S4").Value = 46
$ws.Range("T4").Value = 77
$ws.Range("V4").Value = 78
$ws.Range("W4").Value = 74
$ws.Range("X4").Value = 78
$ws.Range("Y4").Value = 52
$ws.Range("Z4").Value = 79

# Row 5
$ws.Range("P5").Value = 71
$ws.Range("Q5").Value = 66
$ws.Range("R5").Value = 71
$ws.Range("S5").Value = 45
$ws.Range("T5").Value = 77
$ws.Range("V5").Value = 78
$ws.Range("W5").Value = 74
$ws.Range("X5").Value = 78
$ws.Range("Y5").Value = 52
$ws.Range("Z5").Value = 79

# Row 6
$ws.Range("P6").Value = 71
$ws.Range("Q6").Value = 66
$ws.Range("R6").Value = 71
$ws.Range("S6").Value = 38
$ws.Range("T6").Value = 77
$ws.Range("V6").Value = 78
$ws.Range("W6").Value = 73
$ws.Range("X6").Value = 77
$ws.Range("Y6").Value = 44
$ws.Range("Z6").Value = 79

# Row 8
$ws.Range("P8").Value = 100
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = 95
$ws.Range("S8").Value = 50
$ws.Range("T8").Value = 78
$ws.Range("V8").Value = 100
$ws.Range("W8").Value = 4
$ws.Range("X8").Value = 95
$ws.Range("Y8").Value = 50
$ws.Range("Z8").Value = 79

# Row 9
$ws.Range("P9").Value = 72
$ws.Range("Q9").Value = 67
$ws.Range("R9").Value = 72
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 78
$ws.Range("V9").Value = 77
$ws.Range("W9").Value = 72
$ws.Range("X9").Value = 76
$ws.Range("Y9").Value = 6
$ws.Range("Z9").Value = 79

# Row 10
$ws.Range("P10").Value = 71
$ws.Range("Q10").Value = 66
$ws.Range("R10").Value = 72
$ws.Range("S10").Value = 46
$ws.Range("T10").Value = 77
$ws.Range("V10").Value = 77
$ws.Range("W10").Value = 72
$ws.Range("X10").Value = 76
$ws.Range("Y10").Value = 51
$ws.Range("Z10").Value = 79

# Row 11
$ws.Range("P11").Value = 71
$ws.Range("Q11").Value = 67
$ws.Range("R11").Value = 71
$ws.Range("S11").Value = 46
$ws.Range("T11").Value = 77
$ws.Range("V11").Value = 77
$ws.Range("W11").Value = 72
$ws.Range("X11").Value = 76
$ws.Range("Y11").Value = 51
$ws.Range("Z11").Value = 79

# Row 12
$ws.Range("P12").Value = 71
$ws.Range("Q12").Value = 66
$ws.Range("R12").Value = 71
$ws.Range("S12").Value = 39
$ws.Range("T12").Value = 77
$ws.Range("V12").Value = 77
$ws.Range("W12").Value = 72
$ws.Range("X12").Value = 76
$ws.Range("Y12").Value = 44
$ws.Range("Z12").Value = 79

# Row 14
$ws.Range("P14").Value = 99
$ws.Range("Q14").Value = 7
$ws.Range("R14").Value = 90
$ws.Range("S14").Value = 50
$ws.Range("T14").Value = 77
$ws.Range("V14").Value = 100
$ws.Range("W14").Value = 9
$ws.Range("X14").Value = 91
$ws.Range("Y14").Value = 50
$ws.Range("Z14").Value = 78

# Row 15
$ws.Range("P15").Value = 72
$ws.Range("Q15").Value = 66
$ws.Range("R15").Value = 71
$ws.Range("S15").Value = 9
$ws.Range("T15").Value = 77
$ws.Range("V15").Value = 76
$ws.Range("W15").Value = 70
$ws.Range("X15").Value = 75
$ws.Range("Y15").Value = 11
$ws.Range("Z15").Value = 78

# Row 16
$ws.Range("P16").Value = 71
$ws.Range("Q16").Value = 66
$ws.Range("R16").Value = 71
$ws.Range("S16").Value = 47
$ws.Range("T16").Value = 77
$ws.Range("V16").Value = 76
$ws.Range("W16").Value = 71
$ws.Range("X16").Value = 75
$ws.Range("Y16").Value = 51
$ws.Range("Z16").Value = 78

# Row 17
$ws.Range("P17").Value = 71
$ws.Range("Q17").Value = 66
$ws.Range("R17").Value = 71
$ws.Range("S17").Value = 47
$ws.Range("T17").Value = 77
$ws.Range("V17").Value = 76
$ws.Range("W17").Value = 70
$ws.Range("X17").Value = 75
$ws.Range("Y17").Value = 51
$ws.Range("Z17").Value = 78

# Row 18
$ws.Range("P18").Value = 71
$ws.Range("Q18").Value = 66
$ws.Range("R18").Value = 71
$ws.Range("S18").Value = 39
$ws.Range("T18").Value = 76
$ws.Range("V18").Value = 75
$ws.Range("W18").Value = 70
$ws.Range("X18").Value = 74
$ws.Range("Y18").Value = 44
$ws.Range("Z18").Value = 77

# Row 20
$ws.Range("P20").Value = 97
$ws.Range("Q20").Value = 21
$ws.Range("R20").Value = 82
$ws.Range("S20").Value = 50
$ws.Range("T20").Value = 76
$ws.Range("V20").Value = 97
$ws.Range("W20").Value = 22
$ws.Range("X20").Value = 83
$ws.Range("Y20").Value = 50
$ws.Range("Z20").Value = 76

# Row 21
$ws.Range("P21").Value = 71
$ws.Range("Q21").Value = 63
$ws.Range("R21").Value = 70
$ws.Range("S21").Value = 19
$ws.Range("T21").Value = 76
$ws.Range("V21").Value = 76
$ws.Range("W21").Value = 68
$ws.Range("X21").Value = 74
$ws.Range("Y21").Value = 21
$ws.Range("Z21").Value = 76

# Row 22
$ws.Range("P22").Value = 70
$ws.Range("Q22").Value = 63
$ws.Range("R22").Value = 70
$ws.Range("S22").Value = 47
$ws.Range("T22").Value = 75
$ws.Range("V22").Value = 76
$ws.Range("W22").Value = 69
$ws.Range("X22").Value = 73
$ws.Range("Y22").Value = 51
$ws.Range("Z22").Value = 76

# Row 23
$ws.Range("P23").Value = 71
$ws.Range("Q23").Value = 63
$ws.Range("R23").Value = 70
$ws.Range("S23").Value = 47
$ws.Range("T23").Value = 75
$ws.Range("V23").Value = 76
$ws.Range("W23").Value = 68
$ws.Range("X23").Value = 73
$ws.Range("Y23").Value = 51
$ws.Range("Z23").Value = 76

# Row 24
$ws.Range("P24").Value = 71
$ws.Range("Q24").Value = 64
$ws.Range("R24").Value = 70
$ws.Range("S24").Value = 40
$ws.Range("T24").Value = 75
$ws.Range("V24").Value = 75
$ws.Range("W24").Value = 68
$ws.Range("X24").Value = 73
$ws.Range("Y24").Value = 44
$ws.Range("Z24").Value = 76

# Row 26
$ws.Range("P26").Value = 84
$ws.Range("Q26").Value = 48
$ws.Range("R26").Value = 70
$ws.Range("S26").Value = 50
$ws.Range("T26").Value = 74
$ws.Range("V26").Value = 85
$ws.Range("W26").Value = 49
$ws.Range("X26").Value = 70
$ws.Range("Y26").Value = 50
$ws.Range("Z26").Value = 74

# Row 27
$ws.Range("P27").Value = 69
$ws.Range("Q27").Value = 64
$ws.Range("R27").Value = 68
$ws.Range("S27").Value = 38
$ws.Range("T27").Value = 74
$ws.Range("V27").Value = 72
$ws.Range("W27").Value = 67
$ws.Range("X27").Value = 69
$ws.Range("Y27").Value = 41
$ws.Range("Z27").Value = 74

# Row 28
$ws.Range("P28").Value = 69
$ws.Range("Q28").Value = 64
$ws.Range("R28").Value = 68
$ws.Range("S28").Value = 47
$ws.Range("T28").Value = 74
$ws.Range("V28").Value = 72
$ws.Range("W28").Value = 67
$ws.Range("X28").Value = 69
$ws.Range("Y28").Value = 50
$ws.Range("Z28").Value = 74

# Row 29
$ws.Range("P29").Value = 69
$ws.Range("Q29").Value = 64
$ws.Range("R29").Value = 68
$ws.Range("S29").Value = 47
$ws.Range("T29").Value = 74
$ws.Range("V29").Value = 72
$ws.Range("W29").Value = 67
$ws.Range("X29").Value = 69
$ws.Range("Y29").Value = 50
$ws.Range("Z29").Value = 74

# Row 30
$ws.Range("P30").Value = 68
$ws.Range("Q30").Value = 62
$ws.Range("R30").Value = 68
$ws.Range("S30").Value = 40
$ws.Range("T30").Value = 73
$ws.Range("V30").Value = 73
$ws.Range("W30").Value = 67
$ws.Range("X30").Value = 69
$ws.Range("Y30").Value = 43
$ws.Range("Z30").Value = 74

# Row 32
$ws.Range("P32").Value = 72
$ws.Range("Q32").Value = 62
$ws.Range("R32").Value = 67
$ws.Range("S32").Value = 50
$ws.Range("T32").Value = 73
$ws.Range("V32").Value = 73
$ws.Range("W32").Value = 63
$ws.Range("X32").Value = 68
$ws.Range("Y32").Value = 50
$ws.Range("Z32").Value = 74

# Row 33
$ws.Range("P33").Value = 71
$ws.Range("Q33").Value = 62
$ws.Range("R33").Value = 67
$ws.Range("S33").Value = 48
$ws.Range("T33").Value = 73
$ws.Range("V33").Value = 73
$ws.Range("W33").Value = 65
$ws.Range("X33").Value = 68
$ws.Range("Y33").Value = 51
$ws.Range("Z33").Value = 74

# --- View state: select O1 and scroll so column B becomes the first visible
#     column (matches topLeftCell="B1" / selection activeCell="O1") ---
$ws.Range("O1").Select()
$excel.ActiveWindow.ScrollColumn = 2
